# removed "Meinte" from Spoofax interested party
# Slide 9 ("Interested in Integration"), Content Placeholder 2, paragraph 4:
#   "TU Delft Spoofax - Meinte"  ->  "TU Delft Spoofax"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$para = $tr.Paragraphs(4)
$para.Text = "TU Delft Spoofax"
